# Bugfixed evaluation and simulated rt_data for components
# Replace the quarter-label strings in column A (A2:A22) with real dates
# (Dec-31 of each year), formatted as "YYYY-MM-DD HH:MM:SS".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$years = 2004..2024

$range = $ws.Range("A2:A22")

for ($i = 0; $i -lt $years.Count; $i++) {
    $row = 2 + $i
    $cell = $ws.Cells.Item($row, 1)
    $d = Get-Date -Year $years[$i] -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
    $cell.Value2 = [double]$d.ToOADate()
}

$range.NumberFormat = "YYYY-MM-DD HH:MM:SS"
